$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp
$ws.Range('A1').Value = 'Datos actualizados a 5 de Mayo de 2020 a las 01:08'

# Full refreshed & re-sorted (descending by Casos totales) country data table
$data = New-Object 'object[,]' 215,8
$data[0,0] = 'Estados Unidos'
$data[0,1] = 1211213
$data[0,2] = 23091
$data[0,3] = 186991
$data[0,4] = 954632
$data[0,5] = 16043
$data[0,6] = 992
$data[0,7] = 69590
$data[1,0] = 'España'
$data[1,1] = 248301
$data[1,2] = 1179
$data[1,3] = 151633
$data[1,4] = 71240
$data[1,5] = 2254
$data[1,6] = 164
$data[1,7] = 25428
$data[2,0] = 'Italia'
$data[2,1] = 211938
$data[2,2] = 1221
$data[2,3] = 82879
$data[2,4] = 99980
$data[2,5] = 1479
$data[2,6] = 195
$data[2,7] = 29079
$data[3,0] = 'Reino Unido'
$data[3,1] = 190584
$data[3,2] = 3985
$data[3,3] = 0
$data[3,4] = 161506
$data[3,5] = 1559
$data[3,6] = 288
$data[3,7] = 28734
$data[4,0] = 'Francia'
$data[4,1] = 169462
$data[4,2] = 769
$data[4,3] = 51371
$data[4,4] = 92890
$data[4,5] = 3696
$data[4,6] = 306
$data[4,7] = 25201
$data[5,0] = 'Alemania'
$data[5,1] = 166152
$data[5,2] = 488
$data[5,3] = 132700
$data[5,4] = 26459
$data[5,5] = 1949
$data[5,6] = 127
$data[5,7] = 6993
$data[6,0] = 'Rusia'
$data[6,1] = 145268
$data[6,2] = 10581
$data[6,3] = 18095
$data[6,4] = 125817
$data[6,5] = 2300
$data[6,6] = 76
$data[6,7] = 1356
$data[7,0] = 'Turquia'
$data[7,1] = 127659
$data[7,2] = 1614
$data[7,3] = 68166
$data[7,4] = 56032
$data[7,5] = 1384
$data[7,6] = 64
$data[7,7] = 3461
$data[8,0] = 'Brasil'
$data[8,1] = 107844
$data[8,2] = 6697
$data[8,3] = 45815
$data[8,4] = 54701
$data[8,5] = 8318
$data[8,6] = 303
$data[8,7] = 7328
$data[9,0] = 'Iran'
$data[9,1] = 98647
$data[9,2] = 1223
$data[9,3] = 79379
$data[9,4] = 12991
$data[9,5] = 2676
$data[9,6] = 74
$data[9,7] = 6277
$data[10,0] = 'China'
$data[10,1] = 82880
$data[10,2] = 3
$data[10,3] = 77766
$data[10,4] = 481
$data[10,5] = 33
$data[10,6] = 0
$data[10,7] = 4633
$data[11,0] = 'Canada'
$data[11,1] = 60616
$data[11,2] = 1142
$data[11,3] = 25422
$data[11,4] = 31352
$data[11,5] = 557
$data[11,6] = 160
$data[11,7] = 3842
$data[12,0] = 'Belgica'
$data[12,1] = 50267
$data[12,2] = 361
$data[12,3] = 12378
$data[12,4] = 29965
$data[12,5] = 655
$data[12,6] = 80
$data[12,7] = 7924
$data[13,0] = 'Peru'
$data[13,1] = 47372
$data[13,2] = 1444
$data[13,3] = 14427
$data[13,4] = 31601
$data[13,5] = 694
$data[13,6] = 58
$data[13,7] = 1344
$data[14,0] = 'India'
$data[14,1] = 46437
$data[14,2] = 3932
$data[14,3] = 12847
$data[14,4] = 32024
$data[14,5] = 0
$data[14,6] = 175
$data[14,7] = 1566
$data[15,0] = 'Paises Bajos'
$data[15,1] = 40770
$data[15,2] = 199
$data[15,3] = 0
$data[15,4] = 35438
$data[15,5] = 683
$data[15,6] = 26
$data[15,7] = 5082
$data[16,0] = 'Ecuador'
$data[16,1] = 31881
$data[16,2] = 2343
$data[16,3] = 3433
$data[16,4] = 26879
$data[16,5] = 159
$data[16,6] = 5
$data[16,7] = 1569
$data[17,0] = 'Suiza'
$data[17,1] = 29981
$data[17,2] = 76
$data[17,3] = 25200
$data[17,4] = 2997
$data[17,5] = 141
$data[17,6] = 22
$data[17,7] = 1784
$data[18,0] = 'Arabia Saudita'
$data[18,1] = 28656
$data[18,2] = 1645
$data[18,3] = 4476
$data[18,4] = 23989
$data[18,5] = 143
$data[18,6] = 7
$data[18,7] = 191
$data[19,0] = 'Portugal'
$data[19,1] = 25524
$data[19,2] = 242
$data[19,3] = 1712
$data[19,4] = 22749
$data[19,5] = 143
$data[19,6] = 20
$data[19,7] = 1063
$data[20,0] = 'Mexico'
$data[20,1] = 23471
$data[20,2] = 1383
$data[20,3] = 13447
$data[20,4] = 7870
$data[20,5] = 378
$data[20,6] = 93
$data[20,7] = 2154
$data[21,0] = 'Suecia'
$data[21,1] = 22721
$data[21,2] = 404
$data[21,3] = 4074
$data[21,4] = 15878
$data[21,5] = 455
$data[21,6] = 90
$data[21,7] = 2769
$data[22,0] = 'Irlanda'
$data[22,1] = 21772
$data[22,2] = 266
$data[22,3] = 13386
$data[22,4] = 7067
$data[22,5] = 93
$data[22,6] = 16
$data[22,7] = 1319
$data[23,0] = 'Pakistan'
$data[23,1] = 20941
$data[23,2] = 857
$data[23,3] = 5635
$data[23,4] = 14830
$data[23,5] = 111
$data[23,6] = 19
$data[23,7] = 476
$data[24,0] = 'Chile'
$data[24,1] = 20643
$data[24,2] = 980
$data[24,3] = 10415
$data[24,4] = 9958
$data[24,5] = 464
$data[24,6] = 10
$data[24,7] = 270
$data[25,0] = 'Singapur'
$data[25,1] = 18778
$data[25,2] = 573
$data[25,3] = 1457
$data[25,4] = 17303
$data[25,5] = 22
$data[25,6] = 0
$data[25,7] = 18
$data[26,0] = 'Bielorrusia'
$data[26,1] = 17489
$data[26,2] = 784
$data[26,3] = 3259
$data[26,4] = 14127
$data[26,5] = 92
$data[26,6] = 4
$data[26,7] = 103
$data[27,0] = 'Israel'
$data[27,1] = 16246
$data[27,2] = 38
$data[27,3] = 10064
$data[27,4] = 5947
$data[27,5] = 70
$data[27,6] = 3
$data[27,7] = 235
$data[28,0] = 'Catar'
$data[28,1] = 16191
$data[28,2] = 640
$data[28,3] = 1810
$data[28,4] = 14369
$data[28,5] = 72
$data[28,6] = 0
$data[28,7] = 12
$data[29,0] = 'Austria'
$data[29,1] = 15621
$data[29,2] = 24
$data[29,3] = 13316
$data[29,4] = 1705
$data[29,5] = 111
$data[29,6] = 2
$data[29,7] = 600
$data[30,0] = 'Japon'
$data[30,1] = 14877
$data[30,2] = 0
$data[30,3] = 3981
$data[30,4] = 10409
$data[30,5] = 321
$data[30,6] = 0
$data[30,7] = 487
$data[31,0] = 'Emiratos Arabes Unidos'
$data[31,1] = 14730
$data[31,2] = 567
$data[31,3] = 2966
$data[31,4] = 11627
$data[31,5] = 1
$data[31,6] = 11
$data[31,7] = 137
$data[32,0] = 'Polonia'
$data[32,1] = 14006
$data[32,2] = 313
$data[32,3] = 4095
$data[32,4] = 9213
$data[32,5] = 160
$data[32,6] = 20
$data[32,7] = 698
$data[33,0] = 'Rumania'
$data[33,1] = 13512
$data[33,2] = 349
$data[33,3] = 5269
$data[33,4] = 7425
$data[33,5] = 243
$data[33,6] = 28
$data[33,7] = 818
$data[34,0] = 'Ucrania'
$data[34,1] = 12331
$data[34,2] = 418
$data[34,3] = 1619
$data[34,4] = 10409
$data[34,5] = 160
$data[34,6] = 15
$data[34,7] = 303
$data[35,0] = 'Indonesia'
$data[35,1] = 11587
$data[35,2] = 395
$data[35,3] = 1954
$data[35,4] = 8769
$data[35,5] = 0
$data[35,6] = 19
$data[35,7] = 864
$data[36,0] = 'Corea del Sur'
$data[36,1] = 10801
$data[36,2] = 8
$data[36,3] = 9217
$data[36,4] = 1332
$data[36,5] = 55
$data[36,6] = 2
$data[36,7] = 252
$data[37,0] = 'Banglades'
$data[37,1] = 10143
$data[37,2] = 688
$data[37,3] = 1209
$data[37,4] = 8752
$data[37,5] = 1
$data[37,6] = 5
$data[37,7] = 182
$data[38,0] = 'Dinamarca'
$data[38,1] = 9670
$data[38,2] = 147
$data[38,3] = 7088
$data[38,4] = 2089
$data[38,5] = 57
$data[38,6] = 9
$data[38,7] = 493
$data[39,0] = 'Serbia'
$data[39,1] = 9557
$data[39,2] = 93
$data[39,3] = 1574
$data[39,4] = 7786
$data[39,5] = 53
$data[39,6] = 4
$data[39,7] = 197
$data[40,0] = 'Filipinas'
$data[40,1] = 9485
$data[40,2] = 262
$data[40,3] = 1315
$data[40,4] = 7547
$data[40,5] = 31
$data[40,6] = 16
$data[40,7] = 623
$data[41,0] = 'Republica Dominicana'
$data[41,1] = 8235
$data[41,2] = 281
$data[41,3] = 1771
$data[41,4] = 6118
$data[41,5] = 144
$data[41,6] = 13
$data[41,7] = 346
$data[42,0] = 'Colombia'
$data[42,1] = 7973
$data[42,2] = 305
$data[42,3] = 1807
$data[42,4] = 5808
$data[42,5] = 122
$data[42,6] = 18
$data[42,7] = 358
$data[43,0] = 'Noruega'
$data[43,1] = 7884
$data[43,2] = 37
$data[43,3] = 32
$data[43,4] = 7638
$data[43,5] = 27
$data[43,6] = 3
$data[43,7] = 214
$data[44,0] = 'Chequia'
$data[44,1] = 7819
$data[44,2] = 38
$data[44,3] = 3807
$data[44,4] = 3760
$data[44,5] = 58
$data[44,6] = 4
$data[44,7] = 252
$data[45,0] = 'Sudafrica'
$data[45,1] = 7220
$data[45,2] = 437
$data[45,3] = 2746
$data[45,4] = 4336
$data[45,5] = 36
$data[45,6] = 7
$data[45,7] = 138
$data[46,0] = 'Panama'
$data[46,1] = 7197
$data[46,2] = 107
$data[46,3] = 641
$data[46,4] = 6356
$data[46,5] = 91
$data[46,6] = 3
$data[46,7] = 200
$data[47,0] = 'Australia'
$data[47,1] = 6825
$data[47,2] = 24
$data[47,3] = 5859
$data[47,4] = 871
$data[47,5] = 28
$data[47,6] = 0
$data[47,7] = 95
$data[48,0] = 'Egipto'
$data[48,1] = 6813
$data[48,2] = 348
$data[48,3] = 1632
$data[48,4] = 4745
$data[48,5] = 0
$data[48,6] = 7
$data[48,7] = 436
$data[49,0] = 'Malasia'
$data[49,1] = 6353
$data[49,2] = 55
$data[49,3] = 4484
$data[49,4] = 1764
$data[49,5] = 28
$data[49,6] = 0
$data[49,7] = 105
$data[50,0] = 'Finlandia'
$data[50,1] = 5327
$data[50,2] = 73
$data[50,3] = 3500
$data[50,4] = 1587
$data[50,5] = 49
$data[50,6] = 10
$data[50,7] = 240
$data[51,0] = 'Kuwait'
$data[51,1] = 5278
$data[51,2] = 295
$data[51,3] = 1947
$data[51,4] = 3291
$data[51,5] = 79
$data[51,6] = 2
$data[51,7] = 40
$data[52,0] = 'Marruecos'
$data[52,1] = 5053
$data[52,2] = 150
$data[52,3] = 1653
$data[52,4] = 3221
$data[52,5] = 1
$data[52,6] = 5
$data[52,7] = 179
$data[53,0] = 'Argentina'
$data[53,1] = 4783
$data[53,2] = 0
$data[53,3] = 1442
$data[53,4] = 3092
$data[53,5] = 157
$data[53,6] = 3
$data[53,7] = 249
$data[54,0] = 'Argelia'
$data[54,1] = 4648
$data[54,2] = 174
$data[54,3] = 1998
$data[54,4] = 2185
$data[54,5] = 22
$data[54,6] = 2
$data[54,7] = 465
$data[55,0] = 'Moldavia'
$data[55,1] = 4248
$data[55,2] = 127
$data[55,3] = 1423
$data[55,4] = 2693
$data[55,5] = 237
$data[55,6] = 7
$data[55,7] = 132
$data[56,0] = 'Kazajistan'
$data[56,1] = 4049
$data[56,2] = 129
$data[56,3] = 1173
$data[56,4] = 2847
$data[56,5] = 40
$data[56,6] = 2
$data[56,7] = 29
$data[57,0] = 'Luxemburgo'
$data[57,1] = 3828
$data[57,2] = 4
$data[57,3] = 3405
$data[57,4] = 327
$data[57,5] = 21
$data[57,6] = 0
$data[57,7] = 96
$data[58,0] = 'Barein'
$data[58,1] = 3533
$data[58,2] = 150
$data[58,3] = 1744
$data[58,4] = 1781
$data[58,5] = 1
$data[58,6] = 0
$data[58,7] = 8
$data[59,0] = 'Hungria'
$data[59,1] = 3035
$data[59,2] = 37
$data[59,3] = 630
$data[59,4] = 2054
$data[59,5] = 55
$data[59,6] = 11
$data[59,7] = 351
$data[60,0] = 'Tailandia'
$data[60,1] = 2987
$data[60,2] = 18
$data[60,3] = 2740
$data[60,4] = 193
$data[60,5] = 61
$data[60,6] = 0
$data[60,7] = 54
$data[61,0] = 'Afganistan'
$data[61,1] = 2894
$data[61,2] = 190
$data[61,3] = 397
$data[61,4] = 2407
$data[61,5] = 7
$data[61,6] = 5
$data[61,7] = 90
$data[62,0] = 'Ghana'
$data[62,1] = 2719
$data[62,2] = 550
$data[62,3] = 294
$data[62,4] = 2407
$data[62,5] = 4
$data[62,6] = 0
$data[62,7] = 18
$data[63,0] = 'Oman'
$data[63,1] = 2637
$data[63,2] = 69
$data[63,3] = 816
$data[63,4] = 1809
$data[63,5] = 17
$data[63,6] = 0
$data[63,7] = 12
$data[64,0] = 'Grecia'
$data[64,1] = 2632
$data[64,2] = 6
$data[64,3] = 1374
$data[64,4] = 1112
$data[64,5] = 35
$data[64,6] = 2
$data[64,7] = 146
$data[65,0] = 'Nigeria'
$data[65,1] = 2558
$data[65,2] = 0
$data[65,3] = 400
$data[65,4] = 2071
$data[65,5] = 4
$data[65,6] = 0
$data[65,7] = 87
$data[66,0] = 'Armenia'
$data[66,1] = 2507
$data[66,2] = 121
$data[66,3] = 1071
$data[66,4] = 1397
$data[66,5] = 10
$data[66,6] = 4
$data[66,7] = 39
$data[67,0] = 'Irak'
$data[67,1] = 2346
$data[67,2] = 50
$data[67,3] = 1544
$data[67,4] = 704
$data[67,5] = 0
$data[67,6] = 1
$data[67,7] = 98
$data[68,0] = 'Uzbekistan'
$data[68,1] = 2189
$data[68,2] = 40
$data[68,3] = 1405
$data[68,4] = 774
$data[68,5] = 8
$data[68,6] = 0
$data[68,7] = 10
$data[69,0] = 'Croacia'
$data[69,1] = 2101
$data[69,2] = 5
$data[69,3] = 1522
$data[69,4] = 499
$data[69,5] = 15
$data[69,6] = 1
$data[69,7] = 80
$data[70,0] = 'Camerun'
$data[70,1] = 2077
$data[70,2] = 0
$data[70,3] = 953
$data[70,4] = 1060
$data[70,5] = 12
$data[70,6] = 0
$data[70,7] = 64
$data[71,0] = 'Azerbaiyan'
$data[71,1] = 1984
$data[71,2] = 52
$data[71,3] = 1480
$data[71,4] = 478
$data[71,5] = 18
$data[71,6] = 1
$data[71,7] = 26
$data[72,0] = 'Bosnia y Herzegovina'
$data[72,1] = 1926
$data[72,2] = 69
$data[72,3] = 855
$data[72,4] = 993
$data[72,5] = 4
$data[72,6] = 1
$data[72,7] = 78
$data[73,0] = 'Islandia'
$data[73,1] = 1799
$data[73,2] = 0
$data[73,3] = 1723
$data[73,4] = 66
$data[73,5] = 0
$data[73,6] = 0
$data[73,7] = 10
$data[74,0] = 'Guinea'
$data[74,1] = 1710
$data[74,2] = 124
$data[74,3] = 450
$data[74,4] = 1251
$data[74,5] = 0
$data[74,6] = 2
$data[74,7] = 9
$data[75,0] = 'Estonia'
$data[75,1] = 1703
$data[75,2] = 3
$data[75,3] = 259
$data[75,4] = 1389
$data[75,5] = 6
$data[75,6] = 0
$data[75,7] = 55
$data[76,0] = 'Cuba'
$data[76,1] = 1668
$data[76,2] = 19
$data[76,3] = 876
$data[76,4] = 723
$data[76,5] = 9
$data[76,6] = 2
$data[76,7] = 69
$data[77,0] = 'Bulgaria'
$data[77,1] = 1652
$data[77,2] = 34
$data[77,3] = 321
$data[77,4] = 1253
$data[77,5] = 37
$data[77,6] = 5
$data[77,7] = 78
$data[78,0] = 'Bolivia'
$data[78,1] = 1594
$data[78,2] = 124
$data[78,3] = 166
$data[78,4] = 1352
$data[78,5] = 3
$data[78,6] = 5
$data[78,7] = 76
$data[79,0] = 'Republica de Macedonia'
$data[79,1] = 1518
$data[79,2] = 7
$data[79,3] = 992
$data[79,4] = 441
$data[79,5] = 21
$data[79,6] = 1
$data[79,7] = 85
$data[80,0] = 'Nueva Zelanda'
$data[80,1] = 1487
$data[80,2] = 0
$data[80,3] = 1276
$data[80,4] = 191
$data[80,5] = 0
$data[80,6] = 0
$data[80,7] = 20
$data[81,0] = 'Eslovenia'
$data[81,1] = 1439
$data[81,2] = 0
$data[81,3] = 241
$data[81,4] = 1101
$data[81,5] = 20
$data[81,6] = 1
$data[81,7] = 97
$data[82,0] = 'Costa de Marfil'
$data[82,1] = 1432
$data[82,2] = 34
$data[82,3] = 693
$data[82,4] = 722
$data[82,5] = 0
$data[82,6] = 0
$data[82,7] = 17
$data[83,0] = 'Lituania'
$data[83,1] = 1419
$data[83,2] = 9
$data[83,3] = 638
$data[83,4] = 735
$data[83,5] = 17
$data[83,6] = 0
$data[83,7] = 46
$data[84,0] = 'Eslovaquia'
$data[84,1] = 1413
$data[84,2] = 5
$data[84,3] = 643
$data[84,4] = 745
$data[84,5] = 7
$data[84,6] = 1
$data[84,7] = 25
$data[85,0] = 'Senegal'
$data[85,1] = 1271
$data[85,2] = 89
$data[85,3] = 415
$data[85,4] = 846
$data[85,5] = 6
$data[85,6] = 1
$data[85,7] = 10
$data[86,0] = 'Republica de Yibuti'
$data[86,1] = 1116
$data[86,2] = 4
$data[86,3] = 713
$data[86,4] = 401
$data[86,5] = 0
$data[86,6] = 0
$data[86,7] = 2
$data[87,0] = 'Honduras'
$data[87,1] = 1055
$data[87,2] = 45
$data[87,3] = 118
$data[87,4] = 855
$data[87,5] = 10
$data[87,6] = 6
$data[87,7] = 82
$data[88,0] = 'Hong Kong'
$data[88,1] = 1041
$data[88,2] = 1
$data[88,3] = 900
$data[88,4] = 137
$data[88,5] = 1
$data[88,6] = 0
$data[88,7] = 4
$data[89,0] = 'Tunez'
$data[89,1] = 1018
$data[89,2] = 5
$data[89,3] = 406
$data[89,4] = 569
$data[89,5] = 18
$data[89,6] = 1
$data[89,7] = 43
$data[90,0] = 'Letonia'
$data[90,1] = 896
$data[90,2] = 17
$data[90,3] = 348
$data[90,4] = 532
$data[90,5] = 4
$data[90,6] = 0
$data[90,7] = 16
$data[91,0] = 'Republica de Chipre'
$data[91,1] = 874
$data[91,2] = 2
$data[91,3] = 296
$data[91,4] = 563
$data[91,5] = 15
$data[91,6] = 0
$data[91,7] = 15
$data[92,0] = 'Kirguistan'
$data[92,1] = 830
$data[92,2] = 35
$data[92,3] = 575
$data[92,4] = 245
$data[92,5] = 13
$data[92,6] = 0
$data[92,7] = 10
$data[93,0] = 'Albania'
$data[93,1] = 803
$data[93,2] = 8
$data[93,3] = 543
$data[93,4] = 229
$data[93,5] = 7
$data[93,6] = 0
$data[93,7] = 31
$data[94,0] = 'Somalia'
$data[94,1] = 756
$data[94,2] = 34
$data[94,3] = 61
$data[94,4] = 660
$data[94,5] = 2
$data[94,6] = 3
$data[94,7] = 35
$data[95,0] = 'Niger'
$data[95,1] = 755
$data[95,2] = 5
$data[95,3] = 534
$data[95,4] = 184
$data[95,5] = 0
$data[95,6] = 1
$data[95,7] = 37
$data[96,0] = 'Sri Lanka'
$data[96,1] = 751
$data[96,2] = 33
$data[96,3] = 194
$data[96,4] = 549
$data[96,5] = 1
$data[96,6] = 1
$data[96,7] = 8
$data[97,0] = 'Principado de Andorra'
$data[97,1] = 750
$data[97,2] = 2
$data[97,3] = 499
$data[97,4] = 206
$data[97,5] = 16
$data[97,6] = 0
$data[97,7] = 45
$data[98,0] = 'Costa Rica'
$data[98,1] = 742
$data[98,2] = 3
$data[98,3] = 399
$data[98,4] = 337
$data[98,5] = 5
$data[98,6] = 0
$data[98,7] = 6
$data[99,0] = 'Libano'
$data[99,1] = 740
$data[99,2] = 3
$data[99,3] = 200
$data[99,4] = 515
$data[99,5] = 43
$data[99,6] = 0
$data[99,7] = 25
$data[100,0] = 'Crucero'
$data[100,1] = 712
$data[100,2] = 0
$data[100,3] = 645
$data[100,4] = 54
$data[100,5] = 4
$data[100,6] = 0
$data[100,7] = 13
$data[101,0] = 'Guatemala'
$data[101,1] = 703
$data[101,2] = 15
$data[101,3] = 72
$data[101,4] = 614
$data[101,5] = 5
$data[101,6] = 0
$data[101,7] = 17
$data[102,0] = 'Mayotte'
$data[102,1] = 686
$data[102,2] = 36
$data[102,3] = 352
$data[102,4] = 328
$data[102,5] = 6
$data[102,6] = 0
$data[102,7] = 6
$data[103,0] = 'Consejo Danes para los Refugiados'
$data[103,1] = 682
$data[103,2] = 8
$data[103,3] = 80
$data[103,4] = 568
$data[103,5] = 0
$data[103,6] = 1
$data[103,7] = 34
$data[104,0] = 'Sudan'
$data[104,1] = 678
$data[104,2] = 86
$data[104,3] = 61
$data[104,4] = 576
$data[104,5] = 0
$data[104,6] = 0
$data[104,7] = 41
$data[105,0] = 'Burkina Faso'
$data[105,1] = 672
$data[105,2] = 10
$data[105,3] = 545
$data[105,4] = 81
$data[105,5] = 0
$data[105,6] = 1
$data[105,7] = 46
$data[106,0] = 'Uruguay'
$data[106,1] = 655
$data[106,2] = 0
$data[106,3] = 442
$data[106,4] = 196
$data[106,5] = 10
$data[106,6] = 0
$data[106,7] = 17
$data[107,0] = 'Georgia'
$data[107,1] = 593
$data[107,2] = 4
$data[107,3] = 223
$data[107,4] = 361
$data[107,5] = 6
$data[107,6] = 0
$data[107,7] = 9
$data[108,0] = 'San Marino'
$data[108,1] = 582
$data[108,2] = 0
$data[108,3] = 86
$data[108,4] = 455
$data[108,5] = 5
$data[108,6] = 0
$data[108,7] = 41
$data[109,0] = 'Mali'
$data[109,1] = 580
$data[109,2] = 17
$data[109,3] = 223
$data[109,4] = 328
$data[109,5] = 0
$data[109,6] = 2
$data[109,7] = 29
$data[110,0] = 'El Salvador'
$data[110,1] = 555
$data[110,2] = 65
$data[110,3] = 180
$data[110,4] = 362
$data[110,5] = 3
$data[110,6] = 2
$data[110,7] = 13
$data[111,0] = 'Maldivas'
$data[111,1] = 541
$data[111,2] = 14
$data[111,3] = 18
$data[111,4] = 522
$data[111,5] = 2
$data[111,6] = 0
$data[111,7] = 1
$data[112,0] = 'Kenia'
$data[112,1] = 490
$data[112,2] = 25
$data[112,3] = 173
$data[112,4] = 293
$data[112,5] = 2
$data[112,6] = 0
$data[112,7] = 24
$data[113,0] = 'Tanzania'
$data[113,1] = 480
$data[113,2] = 0
$data[113,3] = 167
$data[113,4] = 297
$data[113,5] = 7
$data[113,6] = 0
$data[113,7] = 16
$data[114,0] = 'Malta'
$data[114,1] = 480
$data[114,2] = 3
$data[114,3] = 399
$data[114,4] = 77
$data[114,5] = 1
$data[114,6] = 0
$data[114,7] = 4
$data[115,0] = 'Jamaica'
$data[115,1] = 469
$data[115,2] = 6
$data[115,3] = 38
$data[115,4] = 422
$data[115,5] = 2
$data[115,6] = 1
$data[115,7] = 9
$data[116,0] = 'Jordania'
$data[116,1] = 465
$data[116,2] = 4
$data[116,3] = 370
$data[116,4] = 86
$data[116,5] = 5
$data[116,6] = 0
$data[116,7] = 9
$data[117,0] = 'Taiwan'
$data[117,1] = 438
$data[117,2] = 6
$data[117,3] = 334
$data[117,4] = 98
$data[117,5] = 0
$data[117,6] = 0
$data[117,7] = 6
$data[118,0] = 'Reunion'
$data[118,1] = 424
$data[118,2] = 1
$data[118,3] = 300
$data[118,4] = 124
$data[118,5] = 2
$data[118,6] = 0
$data[118,7] = 0
$data[119,0] = 'Guinea-Bisau'
$data[119,1] = 413
$data[119,2] = 156
$data[119,3] = 19
$data[119,4] = 393
$data[119,5] = 0
$data[119,6] = 0
$data[119,7] = 1
$data[120,0] = 'Paraguay'
$data[120,1] = 396
$data[120,2] = 26
$data[120,3] = 126
$data[120,4] = 260
$data[120,5] = 7
$data[120,6] = 0
$data[120,7] = 10
$data[121,0] = 'Gabon'
$data[121,1] = 367
$data[121,2] = 32
$data[121,3] = 93
$data[121,4] = 268
$data[121,5] = 1
$data[121,6] = 1
$data[121,7] = 6
$data[122,0] = 'Estado de Palestina'
$data[122,1] = 362
$data[122,2] = 9
$data[122,3] = 102
$data[122,4] = 258
$data[122,5] = 0
$data[122,6] = 0
$data[122,7] = 2
$data[123,0] = 'Venezuela'
$data[123,1] = 357
$data[123,2] = 0
$data[123,3] = 158
$data[123,4] = 189
$data[123,5] = 1
$data[123,6] = 0
$data[123,7] = 10
$data[124,0] = 'Mauricio'
$data[124,1] = 332
$data[124,2] = 0
$data[124,3] = 316
$data[124,4] = 6
$data[124,5] = 3
$data[124,6] = 0
$data[124,7] = 10
$data[125,0] = 'Isla de Man'
$data[125,1] = 325
$data[125,2] = 4
$data[125,3] = 271
$data[125,4] = 31
$data[125,5] = 21
$data[125,6] = 1
$data[125,7] = 23
$data[126,0] = 'Montenegro'
$data[126,1] = 323
$data[126,2] = 1
$data[126,3] = 253
$data[126,4] = 62
$data[126,5] = 2
$data[126,6] = 0
$data[126,7] = 8
$data[127,0] = 'Guinea Ecuatorial'
$data[127,1] = 315
$data[127,2] = 0
$data[127,3] = 13
$data[127,4] = 299
$data[127,5] = 0
$data[127,6] = 2
$data[127,7] = 3
$data[128,0] = 'Vietnam'
$data[128,1] = 271
$data[128,2] = 0
$data[128,3] = 219
$data[128,4] = 52
$data[128,5] = 8
$data[128,6] = 0
$data[128,7] = 0
$data[129,0] = 'Ruanda'
$data[129,1] = 261
$data[129,2] = 2
$data[129,3] = 128
$data[129,4] = 133
$data[129,5] = 0
$data[129,6] = 0
$data[129,7] = 0
$data[130,0] = 'Congo'
$data[130,1] = 236
$data[130,2] = 7
$data[130,3] = 26
$data[130,4] = 200
$data[130,5] = 0
$data[130,6] = 1
$data[130,7] = 10
$data[131,0] = 'Tayikistan'
$data[131,1] = 230
$data[131,2] = 102
$data[131,3] = 0
$data[131,4] = 227
$data[131,5] = 0
$data[131,6] = 1
$data[131,7] = 3
$data[132,0] = 'Islas Feroe'
$data[132,1] = 187
$data[132,2] = 0
$data[132,3] = 185
$data[132,4] = 2
$data[132,5] = 0
$data[132,6] = 0
$data[132,7] = 0
$data[133,0] = 'Martinica'
$data[133,1] = 181
$data[133,2] = 2
$data[133,3] = 83
$data[133,4] = 84
$data[133,5] = 5
$data[133,6] = 0
$data[133,7] = 14
$data[134,0] = 'Sierra Leona'
$data[134,1] = 178
$data[134,2] = 12
$data[134,3] = 37
$data[134,4] = 132
$data[134,5] = 0
$data[134,6] = 1
$data[134,7] = 9
$data[135,0] = 'Cabo Verde'
$data[135,1] = 175
$data[135,2] = 10
$data[135,3] = 37
$data[135,4] = 136
$data[135,5] = 0
$data[135,6] = 0
$data[135,7] = 2
$data[136,0] = 'Liberia'
$data[136,1] = 166
$data[136,2] = 8
$data[136,3] = 58
$data[136,4] = 90
$data[136,5] = 0
$data[136,6] = 0
$data[136,7] = 18
$data[137,0] = 'Birmania'
$data[137,1] = 161
$data[137,2] = 6
$data[137,3] = 49
$data[137,4] = 106
$data[137,5] = 0
$data[137,6] = 0
$data[137,7] = 6
$data[138,0] = 'Guadalupe'
$data[138,1] = 152
$data[138,2] = 0
$data[138,3] = 98
$data[138,4] = 42
$data[138,5] = 5
$data[138,6] = 0
$data[138,7] = 12
$data[139,0] = 'Madagascar'
$data[139,1] = 149
$data[139,2] = 0
$data[139,3] = 99
$data[139,4] = 50
$data[139,5] = 1
$data[139,6] = 0
$data[139,7] = 0
$data[140,0] = 'Gibraltar'
$data[140,1] = 144
$data[140,2] = 0
$data[140,3] = 133
$data[140,4] = 11
$data[140,5] = 0
$data[140,6] = 0
$data[140,7] = 0
$data[141,0] = 'Etiopia'
$data[141,1] = 140
$data[141,2] = 5
$data[141,3] = 75
$data[141,4] = 62
$data[141,5] = 0
$data[141,6] = 0
$data[141,7] = 3
$data[142,0] = 'Brunei'
$data[142,1] = 138
$data[142,2] = 0
$data[142,3] = 130
$data[142,4] = 7
$data[142,5] = 2
$data[142,6] = 0
$data[142,7] = 1
$data[143,0] = 'Zambia'
$data[143,1] = 137
$data[143,2] = 13
$data[143,3] = 78
$data[143,4] = 56
$data[143,5] = 1
$data[143,6] = 0
$data[143,7] = 3
$data[144,0] = 'Guayana Francesa'
$data[144,1] = 133
$data[144,2] = 5
$data[144,3] = 100
$data[144,4] = 32
$data[144,5] = 2
$data[144,6] = 0
$data[144,7] = 1
$data[145,0] = 'Togo'
$data[145,1] = 126
$data[145,2] = 2
$data[145,3] = 74
$data[145,4] = 43
$data[145,5] = 0
$data[145,6] = 0
$data[145,7] = 9
$data[146,0] = 'Camboya'
$data[146,1] = 122
$data[146,2] = 0
$data[146,3] = 120
$data[146,4] = 2
$data[146,5] = 1
$data[146,6] = 0
$data[146,7] = 0
$data[147,0] = 'Republica del Chad'
$data[147,1] = 117
$data[147,2] = 0
$data[147,3] = 39
$data[147,4] = 68
$data[147,5] = 0
$data[147,6] = 0
$data[147,7] = 10
$data[148,0] = 'Suazilandia'
$data[148,1] = 116
$data[148,2] = 4
$data[148,3] = 12
$data[148,4] = 103
$data[148,5] = 0
$data[148,6] = 0
$data[148,7] = 1
$data[149,0] = 'Trinidad yTobago'
$data[149,1] = 116
$data[149,2] = 0
$data[149,3] = 99
$data[149,4] = 9
$data[149,5] = 0
$data[149,6] = 0
$data[149,7] = 8
$data[150,0] = 'Bermudas'
$data[150,1] = 115
$data[150,2] = 0
$data[150,3] = 54
$data[150,4] = 54
$data[150,5] = 4
$data[150,6] = 0
$data[150,7] = 7
$data[151,0] = 'Haiti'
$data[151,1] = 100
$data[151,2] = 15
$data[151,3] = 10
$data[151,4] = 79
$data[151,5] = 0
$data[151,6] = 3
$data[151,7] = 11
$data[152,0] = 'Aruba'
$data[152,1] = 100
$data[152,2] = 0
$data[152,3] = 81
$data[152,4] = 17
$data[152,5] = 4
$data[152,6] = 0
$data[152,7] = 2
$data[153,0] = 'Uganda'
$data[153,1] = 97
$data[153,2] = 8
$data[153,3] = 55
$data[153,4] = 42
$data[153,5] = 0
$data[153,6] = 0
$data[153,7] = 0
$data[154,0] = 'Benin'
$data[154,1] = 96
$data[154,2] = 6
$data[154,3] = 50
$data[154,4] = 44
$data[154,5] = 0
$data[154,6] = 0
$data[154,7] = 2
$data[155,0] = 'Monaco'
$data[155,1] = 95
$data[155,2] = 0
$data[155,3] = 78
$data[155,4] = 13
$data[155,5] = 1
$data[155,6] = 0
$data[155,7] = 4
$data[156,0] = 'Guyana'
$data[156,1] = 92
$data[156,2] = 10
$data[156,3] = 27
$data[156,4] = 56
$data[156,5] = 3
$data[156,6] = 0
$data[156,7] = 9
$data[157,0] = 'Republica de Africa Central'
$data[157,1] = 85
$data[157,2] = 13
$data[157,3] = 10
$data[157,4] = 75
$data[157,5] = 0
$data[157,6] = 0
$data[157,7] = 0
$data[158,0] = 'Bahamas'
$data[158,1] = 83
$data[158,2] = 0
$data[158,3] = 25
$data[158,4] = 47
$data[158,5] = 1
$data[158,6] = 0
$data[158,7] = 11
$data[159,0] = 'Barbados'
$data[159,1] = 82
$data[159,2] = 1
$data[159,3] = 46
$data[159,4] = 29
$data[159,5] = 4
$data[159,6] = 0
$data[159,7] = 7
$data[160,0] = 'Liechtenstein'
$data[160,1] = 82
$data[160,2] = 0
$data[160,3] = 55
$data[160,4] = 26
$data[160,5] = 0
$data[160,6] = 0
$data[160,7] = 1
$data[161,0] = 'Mozambique'
$data[161,1] = 80
$data[161,2] = 0
$data[161,3] = 19
$data[161,4] = 61
$data[161,5] = 0
$data[161,6] = 0
$data[161,7] = 0
$data[162,0] = 'San Martin (Parte Holandesa)'
$data[162,1] = 76
$data[162,2] = 0
$data[162,3] = 44
$data[162,4] = 19
$data[162,5] = 7
$data[162,6] = 0
$data[162,7] = 13
$data[163,0] = 'Islas Caimanes'
$data[163,1] = 75
$data[163,2] = 1
$data[163,3] = 14
$data[163,4] = 60
$data[163,5] = 3
$data[163,6] = 0
$data[163,7] = 1
$data[164,0] = 'Nepal'
$data[164,1] = 75
$data[164,2] = 0
$data[164,3] = 16
$data[164,4] = 59
$data[164,5] = 0
$data[164,6] = 0
$data[164,7] = 0
$data[165,0] = 'Libia'
$data[165,1] = 63
$data[165,2] = 0
$data[165,3] = 23
$data[165,4] = 37
$data[165,5] = 0
$data[165,6] = 0
$data[165,7] = 3
$data[166,0] = 'Polinesia Francesa'
$data[166,1] = 58
$data[166,2] = 0
$data[166,3] = 51
$data[166,4] = 7
$data[166,5] = 1
$data[166,6] = 0
$data[166,7] = 0
$data[167,0] = 'Sudan del Sur'
$data[167,1] = 46
$data[167,2] = 0
$data[167,3] = 0
$data[167,4] = 46
$data[167,5] = 0
$data[167,6] = 0
$data[167,7] = 0
$data[168,0] = 'Macao'
$data[168,1] = 45
$data[168,2] = 0
$data[168,3] = 39
$data[168,4] = 6
$data[168,5] = 1
$data[168,6] = 0
$data[168,7] = 0
$data[169,0] = 'Siria'
$data[169,1] = 44
$data[169,2] = 0
$data[169,3] = 27
$data[169,4] = 14
$data[169,5] = 0
$data[169,6] = 0
$data[169,7] = 3
$data[170,0] = 'Malaui'
$data[170,1] = 41
$data[170,2] = 2
$data[170,3] = 9
$data[170,4] = 29
$data[170,5] = 1
$data[170,6] = 0
$data[170,7] = 3
$data[171,0] = 'Mongolia'
$data[171,1] = 40
$data[171,2] = 1
$data[171,3] = 12
$data[171,4] = 28
$data[171,5] = 0
$data[171,6] = 0
$data[171,7] = 0
$data[172,0] = 'Puerto Rico'
$data[172,1] = 39
$data[172,2] = 0
$data[172,3] = 1
$data[172,4] = 36
$data[172,5] = 0
$data[172,6] = 0
$data[172,7] = 2
$data[173,0] = 'Eritrea'
$data[173,1] = 39
$data[173,2] = 0
$data[173,3] = 26
$data[173,4] = 13
$data[173,5] = 0
$data[173,6] = 0
$data[173,7] = 0
$data[174,0] = 'San Martin (Parte Francesa)'
$data[174,1] = 38
$data[174,2] = 0
$data[174,3] = 29
$data[174,4] = 6
$data[174,5] = 1
$data[174,6] = 0
$data[174,7] = 3
$data[175,0] = 'Angola'
$data[175,1] = 35
$data[175,2] = 0
$data[175,3] = 11
$data[175,4] = 22
$data[175,5] = 0
$data[175,6] = 0
$data[175,7] = 2
$data[176,0] = 'Zimbabue'
$data[176,1] = 34
$data[176,2] = 0
$data[176,3] = 5
$data[176,4] = 25
$data[176,5] = 0
$data[176,6] = 0
$data[176,7] = 4
$data[177,0] = 'Guam'
$data[177,1] = 32
$data[177,2] = 0
$data[177,3] = 0
$data[177,4] = 31
$data[177,5] = 0
$data[177,6] = 0
$data[177,7] = 1
$data[178,0] = 'Antigua y Barbuda'
$data[178,1] = 25
$data[178,2] = 0
$data[178,3] = 15
$data[178,4] = 7
$data[178,5] = 1
$data[178,6] = 0
$data[178,7] = 3
$data[179,0] = 'Timor Oriental'
$data[179,1] = 24
$data[179,2] = 0
$data[179,3] = 20
$data[179,4] = 4
$data[179,5] = 0
$data[179,6] = 0
$data[179,7] = 0
$data[180,0] = 'Santo Tome y Principe'
$data[180,1] = 23
$data[180,2] = 7
$data[180,3] = 4
$data[180,4] = 16
$data[180,5] = 0
$data[180,6] = 2
$data[180,7] = 3
$data[181,0] = 'Botsuana'
$data[181,1] = 23
$data[181,2] = 0
$data[181,3] = 8
$data[181,4] = 14
$data[181,5] = 0
$data[181,6] = 0
$data[181,7] = 1
$data[182,0] = 'Granada'
$data[182,1] = 21
$data[182,2] = 0
$data[182,3] = 13
$data[182,4] = 8
$data[182,5] = 4
$data[182,6] = 0
$data[182,7] = 0
$data[183,0] = 'Laos'
$data[183,1] = 19
$data[183,2] = 0
$data[183,3] = 9
$data[183,4] = 10
$data[183,5] = 0
$data[183,6] = 0
$data[183,7] = 0
$data[184,0] = 'Fiyi'
$data[184,1] = 18
$data[184,2] = 0
$data[184,3] = 14
$data[184,4] = 4
$data[184,5] = 0
$data[184,6] = 0
$data[184,7] = 0
$data[185,0] = 'Santa Lucia'
$data[185,1] = 18
$data[185,2] = 0
$data[185,3] = 15
$data[185,4] = 3
$data[185,5] = 0
$data[185,6] = 0
$data[185,7] = 0
$data[186,0] = 'Belice'
$data[186,1] = 18
$data[186,2] = 0
$data[186,3] = 13
$data[186,4] = 3
$data[186,5] = 1
$data[186,6] = 0
$data[186,7] = 2
$data[187,0] = 'Nueva Caledonia'
$data[187,1] = 18
$data[187,2] = 0
$data[187,3] = 17
$data[187,4] = 1
$data[187,5] = 1
$data[187,6] = 0
$data[187,7] = 0
$data[188,0] = 'Islas Virgenes de los Estados Unidos'
$data[188,1] = 17
$data[188,2] = 0
$data[188,3] = 0
$data[188,4] = 17
$data[188,5] = 0
$data[188,6] = 0
$data[188,7] = 0
$data[189,0] = 'San Vicente y las Granadinas'
$data[189,1] = 17
$data[189,2] = 1
$data[189,3] = 9
$data[189,4] = 8
$data[189,5] = 0
$data[189,6] = 0
$data[189,7] = 0
$data[190,0] = 'Gambia'
$data[190,1] = 17
$data[190,2] = 0
$data[190,3] = 9
$data[190,4] = 7
$data[190,5] = 0
$data[190,6] = 0
$data[190,7] = 1
$data[191,0] = 'Namibia'
$data[191,1] = 16
$data[191,2] = 0
$data[191,3] = 8
$data[191,4] = 8
$data[191,5] = 0
$data[191,6] = 0
$data[191,7] = 0
$data[192,0] = 'Dominica'
$data[192,1] = 16
$data[192,2] = 0
$data[192,3] = 13
$data[192,4] = 3
$data[192,5] = 0
$data[192,6] = 0
$data[192,7] = 0
$data[193,0] = 'Curazao'
$data[193,1] = 16
$data[193,2] = 0
$data[193,3] = 13
$data[193,4] = 2
$data[193,5] = 0
$data[193,6] = 0
$data[193,7] = 1
$data[194,0] = 'San Cristobal y Nieves'
$data[194,1] = 15
$data[194,2] = 0
$data[194,3] = 8
$data[194,4] = 7
$data[194,5] = 0
$data[194,6] = 0
$data[194,7] = 0
$data[195,0] = 'Burundi'
$data[195,1] = 15
$data[195,2] = 0
$data[195,3] = 7
$data[195,4] = 7
$data[195,5] = 0
$data[195,6] = 0
$data[195,7] = 1
$data[196,0] = 'Nicaragua'
$data[196,1] = 15
$data[196,2] = 1
$data[196,3] = 7
$data[196,4] = 3
$data[196,5] = 0
$data[196,6] = 2
$data[196,7] = 5
$data[197,0] = 'Islas Malvinas'
$data[197,1] = 13
$data[197,2] = 0
$data[197,3] = 13
$data[197,4] = 0
$data[197,5] = 0
$data[197,6] = 0
$data[197,7] = 0
$data[198,0] = 'Yemen'
$data[198,1] = 12
$data[198,2] = 2
$data[198,3] = 1
$data[198,4] = 9
$data[198,5] = 0
$data[198,6] = 0
$data[198,7] = 2
$data[199,0] = 'Islas Turcas y Caicos'
$data[199,1] = 12
$data[199,2] = 0
$data[199,3] = 6
$data[199,4] = 5
$data[199,5] = 0
$data[199,6] = 0
$data[199,7] = 1
$data[200,0] = 'Santa Sede'
$data[200,1] = 11
$data[200,2] = 0
$data[200,3] = 2
$data[200,4] = 9
$data[200,5] = 0
$data[200,6] = 0
$data[200,7] = 0
$data[201,0] = 'Seychelles'
$data[201,1] = 11
$data[201,2] = 0
$data[201,3] = 6
$data[201,4] = 5
$data[201,5] = 0
$data[201,6] = 0
$data[201,7] = 0
$data[202,0] = 'Montserrat'
$data[202,1] = 11
$data[202,2] = 0
$data[202,3] = 7
$data[202,4] = 3
$data[202,5] = 1
$data[202,6] = 0
$data[202,7] = 1
$data[203,0] = 'Groenlandia'
$data[203,1] = 11
$data[203,2] = 0
$data[203,3] = 11
$data[203,4] = 0
$data[203,5] = 0
$data[203,6] = 0
$data[203,7] = 0
$data[204,0] = 'Surinam'
$data[204,1] = 10
$data[204,2] = 0
$data[204,3] = 9
$data[204,4] = 0
$data[204,5] = 0
$data[204,6] = 0
$data[204,7] = 1
$data[205,0] = 'Papua Nueva Guinea'
$data[205,1] = 8
$data[205,2] = 0
$data[205,3] = 6
$data[205,4] = 2
$data[205,5] = 0
$data[205,6] = 0
$data[205,7] = 0
$data[206,0] = 'Mauritania'
$data[206,1] = 8
$data[206,2] = 0
$data[206,3] = 6
$data[206,4] = 1
$data[206,5] = 0
$data[206,6] = 0
$data[206,7] = 1
$data[207,0] = 'Butan'
$data[207,1] = 7
$data[207,2] = 0
$data[207,3] = 5
$data[207,4] = 2
$data[207,5] = 0
$data[207,6] = 0
$data[207,7] = 0
$data[208,0] = 'Bonaire, San Eustaquio y Saba'
$data[208,1] = 6
$data[208,2] = 0
$data[208,3] = 0
$data[208,4] = 6
$data[208,5] = 0
$data[208,6] = 0
$data[208,7] = 0
$data[209,0] = 'Islas Virgenes Britanicas'
$data[209,1] = 6
$data[209,2] = 0
$data[209,3] = 3
$data[209,4] = 2
$data[209,5] = 0
$data[209,6] = 0
$data[209,7] = 1
$data[210,0] = 'Sahara Occidental'
$data[210,1] = 6
$data[210,2] = 0
$data[210,3] = 5
$data[210,4] = 1
$data[210,5] = 0
$data[210,6] = 0
$data[210,7] = 0
$data[211,0] = 'San Bartolome'
$data[211,1] = 6
$data[211,2] = 0
$data[211,3] = 6
$data[211,4] = 0
$data[211,5] = 0
$data[211,6] = 0
$data[211,7] = 0
$data[212,0] = 'Comoras'
$data[212,1] = 3
$data[212,2] = 0
$data[212,3] = 0
$data[212,4] = 3
$data[212,5] = 0
$data[212,6] = 0
$data[212,7] = 0
$data[213,0] = 'Anguila'
$data[213,1] = 3
$data[213,2] = 0
$data[213,3] = 3
$data[213,4] = 0
$data[213,5] = 0
$data[213,6] = 0
$data[213,7] = 0
$data[214,0] = 'San Pedro y Miquelon'
$data[214,1] = 1
$data[214,2] = 0
$data[214,3] = 0
$data[214,4] = 1
$data[214,5] = 0
$data[214,6] = 0
$data[214,7] = 0

$ws.Range('A4:H218').Value = $data
